# Fix Slide 50 <=
# The example "for" loop printed on Slide 50 used a strict "<" comparison
# (`for(int i = 1; i < 1000; i++)`), which under-counts by one iteration.
# Update the sample code so the loop condition reads "i <= 1000".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(50)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$needle = " < 1000; "
$idx = $tr.Text.IndexOf($needle)

if ($idx -ge 0) {
    $target = $tr.Characters($idx + 1, 3)
    $target.Text = " <= "
}
